# Automatic update of files.
# The "Förändrad" (Changed) date column (C) for every data row (2-31) is
# bumped forward by one day (serial date 46060 -> 46061, i.e. 2026-02-07 -> 2026-02-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
